$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: update Objetivos B/C text (PT objectives paragraph)
$ws.Cells.Item(10,2).Value = 'Possibilitar ao aluno uma introdução ao projeto de software, através de trabalho em temas a serem desenvolvidos em sala de aula. Através de trabalhos desenvolvidos em grupo, ao final do curso o aluno conseguirá criar programas e executáveis em python, com a possibilidade de disponibilização do código em repositórios online'
$ws.Cells.Item(10,3).Value = 'Possibilitar ao aluno uma introdução ao projeto de software, através de trabalho em temas a serem desenvolvidos em sala de aula. Através de trabalhos desenvolvidos em grupo, ao final do curso o aluno conseguirá criar programas e executáveis em python, com a possibilidade de disponibilização do código em repositórios online'

# Rows 13-14: remove old A-column labels (now plain docente rows)
$ws.Cells.Item(13,1).Clear()
$ws.Cells.Item(14,1).Clear()

# Rows 13-14: docente list entries in B/C
$ws.Cells.Item(13,2).Value = '7290967 - Emerson Gonçalves de Melo'
$ws.Cells.Item(13,3).Value = '7290967 - Emerson Gonçalves de Melo'
$ws.Cells.Item(14,2).Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Cells.Item(14,3).Value = '1176388 - Luiz Tadeu Fernandes Eleno'

# Rows 15-18: update labels and PT/EN content
$ws.Cells.Item(15,1).Value = 'Programa resumido:'
$ws.Cells.Item(16,1).Value = 'Short syllabus:'
$ws.Cells.Item(17,1).Value = 'Programa:'
$ws.Cells.Item(18,1).Value = 'Syllabus:'

# Row 17 B/C: new cells, set formatting before value to match column style
$ws.Cells.Item(17,2).Font.Bold = $false
$ws.Cells.Item(17,2).WrapText = $true
$ws.Cells.Item(17,2).Value = 'Desenvolvimento de Software: principais estratégias e recursos. Sistemas controladores de versão: git e variantes; Repositórios online: github, gitlab e variantes; Ferramentas de automatização: bash, make, etc; Gerenciamento de documentação: códigos autocomentados e criação automática de documentação; Integrated Development Environments (IDEs): VSCode, Atom etc.; Gestão de projetos computacionais: etapas de desenvolvimento e implementação; controle de versão; trabalho em grupo; organograma de criação de software.'
$ws.Cells.Item(17,3).Font.Bold = $false
$ws.Cells.Item(17,3).WrapText = $true
$ws.Cells.Item(17,3).Font.Color = 255
$ws.Cells.Item(17,3).Value = 'Desenvolvimento de Software: principais estratégias e recursos. Sistemas controladores de versão: git e variantes; Repositórios online: github, gitlab e variantes; Ferramentas de automatização: bash, make, etc; Gerenciamento de documentação: códigos autocomentados e criação automática de documentação; Integrated Development Environments (IDEs): VSCode, Atom etc.; Gestão de projetos computacionais: etapas de desenvolvimento e implementação; controle de versão; trabalho em grupo; organograma de criação de software.'

# Rows 15, 16, 18: B/C content updates (cells already formatted)
$ws.Cells.Item(15,2).Value = 'Desenvolvimento de Software; Sistemas controladores de versão; Repositórios online; Ferramentas de automatização; Gerenciamento de documentação; Integrated Development Environments (IDEs); Gestão de projetos computacionais'
$ws.Cells.Item(15,3).Value = 'Desenvolvimento de Software; Sistemas controladores de versão; Repositórios online; Ferramentas de automatização; Gerenciamento de documentação; Integrated Development Environments (IDEs); Gestão de projetos computacionais'
$ws.Cells.Item(16,2).Value = 'Software development; Version controller systems; Online repositories; automation tools; Documentation management; Integrated Development Environments (IDEs); Computer project management'
$ws.Cells.Item(16,3).Value = 'Software development; Version controller systems; Online repositories; automation tools; Documentation management; Integrated Development Environments (IDEs); Computer project management'
$ws.Cells.Item(18,2).Value = 'Software development: main strategies and resources. Version controller systems: git and variants; Online repositories: github, gitlab and variants; Automation tools: bash, make, etc; Documentation management: auto-commented codes and automatic documentation creation; Integrated Development Environments (IDEs): VSCode, Atom etc.; Management of computer projects: development and implementation stages; version control; group work; software creation organization chart.'
$ws.Cells.Item(18,3).Value = 'Software development: main strategies and resources. Version controller systems: git and variants; Online repositories: github, gitlab and variants; Automation tools: bash, make, etc; Documentation management: auto-commented codes and automatic documentation creation; Integrated Development Environments (IDEs): VSCode, Atom etc.; Management of computer projects: development and implementation stages; version control; group work; software creation organization chart.'

# Row 19: label change, clear old B/C content
$ws.Cells.Item(19,1).Value = 'Avaliação:'
$ws.Cells.Item(19,2).Clear()
$ws.Cells.Item(19,3).Clear()

# Rows 20-22: label and content updates (cells already formatted)
$ws.Cells.Item(20,1).Value = 'Método:'
$ws.Cells.Item(20,2).Value = 'Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados.'
$ws.Cells.Item(20,3).Value = 'Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados.'
$ws.Cells.Item(21,1).Value = 'Critério:'
$ws.Cells.Item(21,2).Value = 'Apresentação de monografia e defesa de projeto.'
$ws.Cells.Item(21,3).Value = 'Apresentação de monografia e defesa de projeto.'

# Row 22: label update + new B/C cells
$ws.Cells.Item(22,1).Value = 'Norma de recuperação:'
$ws.Cells.Item(22,2).Font.Bold = $false
$ws.Cells.Item(22,2).WrapText = $true
$ws.Cells.Item(22,2).Value = 'Não haverá exame de recuperação.'
$ws.Cells.Item(22,3).Font.Bold = $false
$ws.Cells.Item(22,3).WrapText = $true
$ws.Cells.Item(22,3).Font.Color = 255
$ws.Cells.Item(22,3).Value = 'Não haverá exame de recuperação.'

# Row 23: new A-column label (new cell) + B/C bibliography text update
$ws.Cells.Item(23,1).Font.Bold = $true
$ws.Cells.Item(23,1).Value = 'Bibliografia:'
$ws.Cells.Item(23,2).Value = 'Lambert, K. A. Fundamentos de Python: estruturas de dados. Cengage, 2ed, 2022.Nilo Ney Coutinho Menezes. Introdução à Programação com Python: Algoritmos e Lógica de Programação Para Iniciantes, 3a ed, 2019.Ramalho, L. Python Fluente. O’Reilly-Novatec, 2015Downey, A. B. Pense em Python. O’Reilly-Novatec, 2016.STEWART, J. M. Python for scientists. Cambridge University Press, 2014.TELLES, M. Python Power, Boston: Thomson Course Technology PTR, 2008.LUTZ, Mark. Programming Python, 3a ed, Sebastopol, CA: O’Reilly Media, 2006.MCGREGGOR, D. M. Mastering matplotlib. Birmingham, UK: Packt Publishing, 2015.'
$ws.Cells.Item(23,3).Value = 'Lambert, K. A. Fundamentos de Python: estruturas de dados. Cengage, 2ed, 2022.Nilo Ney Coutinho Menezes. Introdução à Programação com Python: Algoritmos e Lógica de Programação Para Iniciantes, 3a ed, 2019.Ramalho, L. Python Fluente. O’Reilly-Novatec, 2015Downey, A. B. Pense em Python. O’Reilly-Novatec, 2016.STEWART, J. M. Python for scientists. Cambridge University Press, 2014.TELLES, M. Python Power, Boston: Thomson Course Technology PTR, 2008.LUTZ, Mark. Programming Python, 3a ed, Sebastopol, CA: O’Reilly Media, 2006.MCGREGGOR, D. M. Mastering matplotlib. Birmingham, UK: Packt Publishing, 2015.'

# Row 24: new A-column label (new cell)
$ws.Cells.Item(24,1).Font.Bold = $true
$ws.Cells.Item(24,1).Value = 'Requisitos:'

# Row 25: new B/C cells (requisito text)
$ws.Cells.Item(25,2).Font.Bold = $false
$ws.Cells.Item(25,2).WrapText = $true
$ws.Cells.Item(25,2).Value = 'LOM3260 -  Computação Científica em Python  (Requisito)' + "`n" + ''
$ws.Cells.Item(25,3).Font.Bold = $false
$ws.Cells.Item(25,3).WrapText = $true
$ws.Cells.Item(25,3).Font.Color = 255
$ws.Cells.Item(25,3).Value = 'LOM3260 -  Computação Científica em Python  (Requisito)' + "`n" + ''

# Row heights
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(25).RowHeight = 30
